$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Cells.Item(3, 8).Value = 16157
$ws.Cells.Item(3, 10).Value = 16157
$ws.Cells.Item(3, 12).Value = 16157
$ws.Cells.Item(3, 14).Value = -16385

# Row 11
$ws.Cells.Item(11, 8).Value = 76.15385000000001
$ws.Cells.Item(11, 9).Value = 76.15385000000001
$ws.Cells.Item(11, 11).Value = 76.15385000000001
$ws.Cells.Item(11, 13).Value = 63.84614999999999

# Row 28
$ws.Cells.Item(28, 8).Value = 2946
$ws.Cells.Item(28, 9).Value = 1025.4375
$ws.Cells.Item(28, 11).Value = 1025.4375
$ws.Cells.Item(28, 13).Value = -540.4375

# Row 38
$ws.Cells.Item(38, 8).Value = 2384.7273
$ws.Cells.Item(38, 9).Value = 47.166668
$ws.Cells.Item(38, 10).Value = 5189.8
$ws.Cells.Item(38, 11).Value = 141.500004
$ws.Cells.Item(38, 12).Value = 15569.4
$ws.Cells.Item(38, 13).Value = 230.499996
$ws.Cells.Item(38, 14).Value = -16313.4

# Row 42
$ws.Cells.Item(42, 8).Value = 863.2
$ws.Cells.Item(42, 10).Value = 450
$ws.Cells.Item(42, 12).Value = 1350
$ws.Cells.Item(42, 14).Value = -1810

# Row 68
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()

# Row 71
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()

# Row 80
$ws.Cells.Item(80, 8).Value = 198.63637
$ws.Cells.Item(80, 9).Value = 85
$ws.Cells.Item(80, 10).Value = 397.5
$ws.Cells.Item(80, 11).Value = 255
$ws.Cells.Item(80, 12).Value = 1192.5
$ws.Cells.Item(80, 13).Value = 743
$ws.Cells.Item(80, 14).Value = -3188.5

# Row 83
$ws.Cells.Item(83, 8).Value = 198.63637
$ws.Cells.Item(83, 9).Value = 85
$ws.Cells.Item(83, 10).Value = 397.5
$ws.Cells.Item(83, 11).Value = 765
$ws.Cells.Item(83, 12).Value = 3577.5
$ws.Cells.Item(83, 13).Value = 4227
$ws.Cells.Item(83, 14).Value = -13561.5

# Row 93
$ws.Cells.Item(93, 8).Value = 25601
$ws.Cells.Item(93, 10).Value = 25601
$ws.Cells.Item(93, 12).Value = 25601
$ws.Cells.Item(93, 14).Value = -30593

# Row 102
$ws.Cells.Item(102, 8).Value = 16157
$ws.Cells.Item(102, 10).Value = 16157
$ws.Cells.Item(102, 12).Value = 16157
$ws.Cells.Item(102, 14).Value = -22647

# Row 107
$ws.Cells.Item(107, 8).Value = 1650.3077
$ws.Cells.Item(107, 9).Value = 2150.4443
$ws.Cells.Item(107, 10).Value = 525
$ws.Cells.Item(107, 11).Value = 2150.4443
$ws.Cells.Item(107, 12).Value = 525
$ws.Cells.Item(107, 13).Value = -230.4443000000001
$ws.Cells.Item(107, 14).Value = -4365

# Row 113
$ws.Cells.Item(113, 8).Value = 5621.4287

# Row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 13).ClearContents()

# Row 137
$ws.Cells.Item(137, 8).Value = 1957.4117
$ws.Cells.Item(137, 9).Value = 946.7778
$ws.Cells.Item(137, 10).Value = 3094.375
$ws.Cells.Item(137, 11).Value = 2840.3334
$ws.Cells.Item(137, 12).Value = 9283.125
$ws.Cells.Item(137, 13).Value = -290.3334
$ws.Cells.Item(137, 14).Value = -14383.125

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 1805.5264
$ws.Cells.Item(61, 9).Value = 1284.7858
$ws.Cells.Item(61, 10).Value = 3263.6
$ws.Cells.Item(61, 11).Value = 1284.7858
$ws.Cells.Item(61, 12).Value = 3263.6
$ws.Cells.Item(61, 13).Value = -1072.7858
$ws.Cells.Item(61, 14).Value = -3687.6

# Row 102
$ws.Cells.Item(102, 8).Value = 8250
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 13).ClearContents()

# Row 110
$ws.Cells.Item(110, 8).Value = 2289.9
$ws.Cells.Item(110, 9).Value = 2173.25
$ws.Cells.Item(110, 11).Value = 2173.25
$ws.Cells.Item(110, 13).Value = -128.25

# Row 136
$ws.Cells.Item(136, 8).Value = 1805.5264
$ws.Cells.Item(136, 9).Value = 1284.7858
$ws.Cells.Item(136, 10).Value = 3263.6
$ws.Cells.Item(136, 11).Value = 3854.3574
$ws.Cells.Item(136, 12).Value = 9790.799999999999
$ws.Cells.Item(136, 13).Value = -1304.3574
$ws.Cells.Item(136, 14).Value = -14890.8

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 1328.75
$ws.Cells.Item(99, 9).Value = 1468.5714
$ws.Cells.Item(99, 11).Value = 1468.5714
$ws.Cells.Item(99, 13).Value = 29.42859999999996

# Row 105
$ws.Cells.Item(105, 8).Value = 1273.1111
$ws.Cells.Item(105, 9).Value = 1182.8125
$ws.Cells.Item(105, 10).Value = 1995.5
$ws.Cells.Item(105, 11).Value = 1182.8125
$ws.Cells.Item(105, 12).Value = 1995.5
$ws.Cells.Item(105, 13).Value = 564.1875
$ws.Cells.Item(105, 14).Value = -5489.5

# Row 134
$ws.Cells.Item(134, 8).Value = 2197.9
$ws.Cells.Item(134, 9).Value = 865
$ws.Cells.Item(134, 10).Value = 5308
$ws.Cells.Item(134, 11).Value = 2595
$ws.Cells.Item(134, 12).Value = 15924
$ws.Cells.Item(134, 13).Value = -60
$ws.Cells.Item(134, 14).Value = -20994

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Cells.Item(107, 8).Value = 495.36365
$ws.Cells.Item(107, 9).Value = 238.35294
$ws.Cells.Item(107, 10).Value = 1369.2
$ws.Cells.Item(107, 11).Value = 238.35294
$ws.Cells.Item(107, 12).Value = 1369.2
$ws.Cells.Item(107, 13).Value = 1681.64706
$ws.Cells.Item(107, 14).Value = -5209.2

# Row 132
$ws.Cells.Item(132, 8).Value = 1823.2667
$ws.Cells.Item(132, 9).Value = 1918.6666
$ws.Cells.Item(132, 10).Value = 1441.6666
$ws.Cells.Item(132, 11).Value = 5755.9998
$ws.Cells.Item(132, 12).Value = 4324.9998
$ws.Cells.Item(132, 13).Value = -3225.9998
$ws.Cells.Item(132, 14).Value = -9384.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Cells.Item(34, 8).Value = 1075
$ws.Cells.Item(34, 10).Value = 1342.2
$ws.Cells.Item(34, 12).Value = 4026.6
$ws.Cells.Item(34, 14).Value = -4194.6

# Row 38
$ws.Cells.Item(38, 8).Value = 454.86365
$ws.Cells.Item(38, 10).Value = 537.25
$ws.Cells.Item(38, 12).Value = 1611.75
$ws.Cells.Item(38, 14).Value = -2305.75

# Row 39
$ws.Cells.Item(39, 8).Value = 5346.6665
$ws.Cells.Item(39, 9).Value = 3798
$ws.Cells.Item(39, 10).Value = 5487.4546
$ws.Cells.Item(39, 11).Value = 11394
$ws.Cells.Item(39, 12).Value = 16462.3638
$ws.Cells.Item(39, 13).Value = -11100
$ws.Cells.Item(39, 14).Value = -17050.3638

# Row 55
$ws.Cells.Item(55, 8).Value = 5221.2666
$ws.Cells.Item(55, 9).Value = 4
$ws.Cells.Item(55, 10).Value = 5593.9287
$ws.Cells.Item(55, 11).Value = 12
$ws.Cells.Item(55, 12).Value = 16781.7861
$ws.Cells.Item(55, 13).Value = 165
$ws.Cells.Item(55, 14).Value = -17135.7861

# Row 101
$ws.Cells.Item(101, 8).Value = 12000
$ws.Cells.Item(101, 10).Value = 12000
$ws.Cells.Item(101, 12).Value = 36000
$ws.Cells.Item(101, 14).Value = -40868

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Cells.Item(126, 8).Value = 3129.2
$ws.Cells.Item(126, 9).Value = 2375.2856
$ws.Cells.Item(126, 10).Value = 4888.3335
$ws.Cells.Item(126, 11).Value = 7125.8568
$ws.Cells.Item(126, 12).Value = 14665.0005
$ws.Cells.Item(126, 13).Value = -4655.8568
$ws.Cells.Item(126, 14).Value = -19605.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Cells.Item(20, 8).Value = 50949.95
$ws.Cells.Item(20, 10).Value = 50949.95
$ws.Cells.Item(20, 12).Value = 50949.95
$ws.Cells.Item(20, 14).Value = -51401.95

# Row 25
$ws.Cells.Item(25, 8).Value = 3752.75
$ws.Cells.Item(25, 9).Value = 3752.75
$ws.Cells.Item(25, 11).Value = 3752.75
$ws.Cells.Item(25, 13).Value = -3522.75

# Row 42
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()

# Row 46
$ws.Cells.Item(46, 8).Value = 5832.6665
$ws.Cells.Item(46, 9).Value = 1998
$ws.Cells.Item(46, 10).Value = 7750
$ws.Cells.Item(46, 11).Value = 1998
$ws.Cells.Item(46, 12).Value = 7750
$ws.Cells.Item(46, 13).Value = -1810
$ws.Cells.Item(46, 14).Value = -8126

# Row 49
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()

# Row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()

# Row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()

# Row 132
$ws.Cells.Item(132, 8).Value = 4049.1
$ws.Cells.Item(132, 9).Value = 3957.3333
$ws.Cells.Item(132, 11).Value = 11871.9999
$ws.Cells.Item(132, 13).Value = -9341.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Cells.Item(3, 8).Value = 57004
$ws.Cells.Item(3, 10).Value = 57004
$ws.Cells.Item(3, 12).Value = 57004
$ws.Cells.Item(3, 14).Value = -57232

# Row 132
$ws.Cells.Item(132, 8).Value = 1468.2413
$ws.Cells.Item(132, 9).Value = 1377.8214
$ws.Cells.Item(132, 11).Value = 4133.4642
$ws.Cells.Item(132, 13).Value = -1603.4642

# Row 136
$ws.Cells.Item(136, 8).Value = 3504.0881
$ws.Cells.Item(136, 9).Value = 2754.261
$ws.Cells.Item(136, 11).Value = 8262.782999999999
$ws.Cells.Item(136, 13).Value = -5712.782999999999
